$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 291, shifting existing rows 291-428 down to 292-429.
$ws.Rows.Item(291).Insert()

# Populate the new row 291 with the weekly record.
$ws.Range("A291").Value = 10
$ws.Range("B291").Value = "Vega Modelo de Temuco"
$ws.Range("C291").Value = "La Araucanía"
$ws.Range("D291").Value = 45202
$ws.Range("D291").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E291").Value = 9
$ws.Range("F291").Value = 100112039
$ws.Range("G291").Value = "Ciboulette"
$ws.Range("H291").Value = "Sin especificar"
$ws.Range("I291").Value = "Primera"
$ws.Range("J291").Value = 45
$ws.Range("K291").Value = 5000
$ws.Range("L291").Value = 5000
$ws.Range("M291").Value = 5000
$ws.Range("N291").Value = "$/docena de atados"
$ws.Range("O291").Value = "Provincia de Cautín"
$ws.Range("P291").Value = 1667
$ws.Range("Q291").Value = 3
$ws.Range("R291").Value = "Hortaliza"
